$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resumos")

$ws.Range("F5").Value = 83
$ws.Range("F6").Formula = "=CEILING(85*0.9,1)"
$ws.Range("F7").Value = "-"
$ws.Range("F8").Value = "-"
$ws.Range("F9").Value = 96
$ws.Range("F10").Value = 91
$ws.Range("F11").Value = 86
$ws.Range("F12").Value = "-"
$ws.Range("F13").Value = 86
$ws.Range("F14").Value = "-"
$ws.Range("F15").Value = "-"
$ws.Range("F16").Formula = "=CEILING(57*0.3,1)"
$ws.Range("F17").Value = 97
$ws.Range("F18").Formula = "=CEILING(100*0.9,1)"
$ws.Range("F19").Value = "-"
$ws.Range("F20").Value = 96
$ws.Range("F21").Value = 95
$ws.Range("F22").Value = 81
$ws.Range("F23").Formula = "=CEILING(92*0.9,1)"
$ws.Range("F24").Value = "-"
$ws.Range("F25").Value = 91
$ws.Range("F26").Value = 88
$ws.Range("F27").Value = 87
$ws.Range("F28").Value = 99
$ws.Range("F29").Value = 90
$ws.Range("F30").Value = 94
$ws.Range("F31").Value = 92
$ws.Range("F32").Value = 93
$ws.Range("F33").Value = 83
$ws.Range("F34").Value = 95
$ws.Range("F35").Value = "-"
$ws.Range("F36").Formula = "=CEILING(87*0.9,1)"
$ws.Range("F37").Value = 95
$ws.Range("F38").Value = 94
$ws.Range("F39").Value = 86
$ws.Range("F40").Value = "-"
$ws.Range("F41").Value = 94
$ws.Range("F42").Formula = "=CEILING(85*0.8,1)"
$ws.Range("F43").Value = 93
$ws.Range("F44").Value = 97
$ws.Range("F45").Value = "-"

$ws.Range("D22").Select()
Write-Output "done"
